$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Semiconductors")
$ws2 = $wb.Worksheets.Item("Networking")
$ws3 = $wb.Worksheets.Item("Electronics-Computers")

# Insert 4 blank rows before row 12 on Electronics-Computers (pushes old row12->16, 22-25->26-29)
$ws3.Rows("12:15").Insert()

# Order mirrors the original authoring sequence (matches sharedStrings order)
$ws2.Range("A7").Value = "x"
$ws2.Range("B7").Value = "Arista"
$ws2.Range("C7").Value = "ANET"

$ws1.Range("A26").Value = "x"
$ws1.Range("B26").Value = "Infineon"
$ws1.Range("C26").Value = "IFX GR"

$ws3.Range("A12").Value = "x"
$ws3.Range("B12").Value = "Canon"
$ws3.Range("C12").Value = "7751 JP"

$ws1.Range("A27").Value = "x"
$ws1.Range("B27").Value = "STMicro"
$ws1.Range("C27").Value = "STM"

$ws3.Range("A13").Value = "x"
$ws3.Range("B13").Value = "Foxconn"
$ws3.Range("C13").Value = "601138 CH"

$ws3.Range("A14").Value = "x"
$ws3.Range("B14").Value = "FujiFilm"
$ws3.Range("C14").Value = "4901 JP"

$ws1.Range("A28").Value = "x"
$ws1.Range("B28").Value = "TCL Zhonghuan"
$ws1.Range("C28").Value = "002129 CH"

$ws3.Range("A15").Value = "x"
$ws3.Range("B15").Value = "Fujitsu"
$ws3.Range("C15").Value = "6702 JP"

$ws2.Range("A8").Value = "x"
$ws2.Range("B8").Value = "Nokia"
$ws2.Range("C8").Value = "NOKIA FH"

$ws1.Range("A29").Value = "x"
$ws1.Range("B29").Value = "GlobalFoundries"
$ws1.Range("C29").Value = "GFS"

$ws1.Range("B36").Value = "Private"
$ws1.Range("B37").Value = "Pasqal"
